# Generate Report for Handback
#
# - Flip every "Ready for handoff" status cell to "Handed back: in sync with en-US"
# - Stamp the actual handback datetime (was the zero-date placeholder) into the
#   "Latest Handback DateTime" column for both locale sheets
# - Populate the (until now empty) "Latest Target File" / "Latest Handback File"
#   columns with their hyperlinked file names, now that handback has happened

$wb = $excel.ActiveWorkbook

# ---- 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ----
# This phrase is shared by the Overview sheet (zh-cn/de-de status columns) and by
# the Status column on each per-locale sheet, so just replace it everywhere.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null
}

# ---- 2. Per-locale sheets: handback datetime + target/handback file columns ----
$localeSheets = @(
    @{ Name = "zh-cn"; Xlf = "e1fe6ad5-9980-499d-9b11-6c0e796de5e6.128aa420178ecbb02ec5442d57c71d2f6af2255d.zh-cn.xlf"; HandbackDate = "2016-03-17 03:17:45" },
    @{ Name = "de-de"; Xlf = "e1fe6ad5-9980-499d-9b11-6c0e796de5e6.128aa420178ecbb02ec5442d57c71d2f6af2255d.de-de.xlf"; HandbackDate = "2016-03-17 03:17:58" }
)

$mdName = "e1fe6ad5-9980-499d-9b11-6c0e796de5e6.md"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ca67744c28f809f4c9709c15c545fc8a3cba8427/e2e/$mdName"

foreach ($locale in $localeSheets) {
    $ws = $wb.Worksheets.Item($locale.Name)
    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/db32c318a11fcfdc09a1a5e8e8cf58ad3fc37644/ol-handoff/OpenLocalizationTestOrg/oltest." + $locale.Name + "/xinjiang/ht/" + $locale.Xlf

    foreach ($row in 2, 3) {
        # Latest Target File (F) - points at the same source .md file as column A
        $fCell = $ws.Cells.Item($row, 6)
        $fCell.Value = $mdName
        $ws.Hyperlinks.Add($fCell, $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null

        # Latest Handback File (G) - points at the localized .xlf handback file
        $gCell = $ws.Cells.Item($row, 7)
        $gCell.Value = $locale.Xlf
        $ws.Hyperlinks.Add($gCell, $xlfUrl, [Type]::Missing, [Type]::Missing, $locale.Xlf) | Out-Null

        # Latest Handback DateTime (H) - stamp the real handback time
        $hCell = $ws.Cells.Item($row, 8)
        $hCell.Value = $locale.HandbackDate
    }
}
